# Update "想去人数" (want-to-go count) values in column F across the
# 展览 (Exhibition), 演出 (Performance), 本地生活 (Local Life) sheets, and
# propagate the same updates to the aggregated 全部类型 (All Types) sheet.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsPerformance = $wb.Worksheets.Item("演出")
$wsLocalLife = $wb.Worksheets.Item("本地生活")
$wsAllTypes = $wb.Worksheets.Item("全部类型")

# --- 展览 (Exhibition) sheet updates ---
$wsExhibition.Range("F11").Value = 7166
$wsExhibition.Range("F14").Value = 1186
$wsExhibition.Range("F21").Value = 731
$wsExhibition.Range("F22").Value = 11
$wsExhibition.Range("F23").Value = 45
$wsExhibition.Range("F24").Value = 119
$wsExhibition.Range("F32").Value = 78
$wsExhibition.Range("F33").Value = 2054
$wsExhibition.Range("F38").Value = 555
$wsExhibition.Range("F39").Value = 11

# --- 演出 (Performance) sheet updates ---
$wsPerformance.Range("F2").Value = 661
$wsPerformance.Range("F5").Value = 292
$wsPerformance.Range("F6").Value = 306
$wsPerformance.Range("F8").Value = 50

# --- 本地生活 (Local Life) sheet updates ---
$wsLocalLife.Range("F2").Value = 361

# --- 全部类型 (All Types) sheet updates (aggregated view) ---
$wsAllTypes.Range("F2").Value = 361
$wsAllTypes.Range("F7").Value = 661
$wsAllTypes.Range("F13").Value = 7166
$wsAllTypes.Range("F17").Value = 1186
$wsAllTypes.Range("F25").Value = 292
$wsAllTypes.Range("F26").Value = 306
$wsAllTypes.Range("F28").Value = 731
$wsAllTypes.Range("F29").Value = 11
$wsAllTypes.Range("F30").Value = 45
$wsAllTypes.Range("F31").Value = 119
$wsAllTypes.Range("F32").Value = 50
$wsAllTypes.Range("F42").Value = 78
$wsAllTypes.Range("F43").Value = 2054
$wsAllTypes.Range("F48").Value = 555
$wsAllTypes.Range("F49").Value = 11
